$d = $word.ActiveDocument

# ---- 1. Insert the new scenario tables + user-stories content -----------
# Select the final (empty) trailing paragraph; InsertXML replaces its
# content, so we append a fresh empty <w:p/> at the end of our payload to
# recreate it after our new content.
$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range

$newContentXml = @'
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblInd w:w="0" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="9350"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="E2EFD9" w:themeFill="accent6" w:themeFillTint="33"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
              </w:rPr>
              <w:t>Benjamin, an R&amp;D chemist</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblInd w:w="0" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="9350"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>
            </w:tcBorders>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Benjamin is well versed in the game of football but hasn’t attended live game or watch a game on television since he attended college as an undergrade. Having never played fantasy football but an avid gamer, he wants an app that has a low learning curve and easy to grasp. He plans on using the tutorial feature to get a sense of what it is like to participate in fantasy football league.  Upon using the tutorial feature, a perusal through the FAQ is there to answer any other queries about the general rule and guidelines to fantasy football.</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblInd w:w="0" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="9350"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Benjamin doesn’t want to deal with the overhead of starting and creating a league but would gladly accept an invitation to join a league that has been already created and looking to add new members. He logs onto the application using his credentials and check for the invitations he has received. If his list of invitations is empty and he can have the app the randomly select a league with an open invitation.</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>He decides that he wants to have a more interactive approach to the league he plays in, so he checks to see what other leagues are actively looking to add players. He selects one and confirms his selection and a message is displayed saying that he has successfully been added to a fantasy league and from then on, he’ll get notifications about drafts, deadlines, player trades and scores updates, match-ups from the other members of the league.</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblInd w:w="0" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="9350"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="9350" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
              <w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>
            </w:tcBorders>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Benjamin wants to use this app as a way to get back into the fandom of football and plans to use the in-app ticker for scores updates and headlines of what’s going on in the real football world. At times he doesn’t want to have to manually update his team weekly before each gameday. He selects the team management tab and choose to have the software automatically removed injured players from his roster and replaces them with a comparable player from the bench if possible. He also selects to have the app automatically select a player of the exact position from the player pool, if there isn’t a comparable player from his bench.</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>User Stories:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As an admin, I want to be able to add or remove players from my league.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As an admin, I don’t want members to be able to drop out of the league once the season starts, so that the fixture list remains consistent.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As a general manager, I want to able to have my injured players automatically removed from the starting lineup.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As an admin, I want to be able to send out invitations to other members about joining my league.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As a general manager, I want to be able to participate in different leagues.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As a general manager, I want to be notified when a new member is added to the league.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As a general manager, I want to be notified of my upcoming fixture.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>As a general manager, I want to able to compare the points of starting line-up vs bench so that I can better predict which change is necessary.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p/>
'@

$insertRange.InsertXML($newContentXml)

# ---- 2. Register the "Table Grid" style referenced by the new tables ----
# Build it the same way Word does when a user applies the built-in "Table
# Grid" style from the Table Design gallery: create a scratch table, apply
# the named style (which mints the style definition), tune the few
# properties that differ from Word's out-of-the-box defaults, then discard
# the scratch table (the style persists in the styles part).
$styleProbeRange = $d.Paragraphs.Last.Range
$probeTable = $d.Tables.Add($styleProbeRange, 1, 1)
$probeTable.Style = "Table Grid"

$gridStyle = $d.Styles("Table Grid")
$gridStyle.Priority = 39
$gridStyle.ParagraphFormat.SpaceAfter = 0
$gridStyle.ParagraphFormat.LineSpacingRule = 0

$probeTable.Delete()

Write-Host "Applied scenario tables + user stories edit."
